$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.363.10'
$ws.Range('E2').Value = '  -1.63%  '
$ws.Range('D3').Value = '2.923.96'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '''371.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.11%  '
$ws.Range('D6').Value = '''104.21'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.91%  '
$ws.Range('D7').Value = '''0.543'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.93%  '
$ws.Range('D8').Value = '''1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '''0.592'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.31%  '
$ws.Range('D10').Value = '''37.26'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.17%  '
$ws.Range('E11').Value = '  +1.00%  '
$ws.Range('D12').Value = '''0.0841'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.08%  '
$ws.Range('D13').Value = '''18.53'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.93%  '
$ws.Range('D14').Value = '3.379.93'
$ws.Range('E14').Value = '  -0.50%  '
$ws.Range('D15').Value = '''7.42'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.76%  '
$ws.Range('D16').Value = '2.916.92'
$ws.Range('E16').Value = '  -0.38%  '
$ws.Range('D17').Value = '''0.955'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.31%  '
$ws.Range('D18').Value = '51.290.22'
$ws.Range('E18').Value = '  -1.68%  '
$ws.Range('E19').Value = '  -5.09%  '
$ws.Range('D20').Value = '''7.31'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.42%  '
$ws.Range('D21').Value = '''13.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.41%  '
$ws.Range('D22').Value = '0.0₃0950'
$ws.Range('E22').Value = '  -3.23%  '
$ws.Range('D23').Value = '''68.52'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.61%  '
$ws.Range('D24').Value = '''261.42'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.25%  '
$ws.Range('D25').Value = '''2.71'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.70%  '
$ws.Range('E26').Value = '  +4.37%  '
$ws.Range('E27').Value = '  -2.58%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('D29').Value = '''7.45'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.36%  '
$ws.Range('D30').Value = '''26.06'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.14%  '
$ws.Range('E31').Value = '  -2.11%  '
$ws.Range('D32').Value = '''6.16'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.47%  '
$ws.Range('D33').Value = '''9.97'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.66%  '
$ws.Range('D34').Value = '''35.42'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.98%  '
$ws.Range('E35').Value = '  -5.74%  '
$ws.Range('D36').Value = '''50.78'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.28%  '
$ws.Range('E37').Value = '  +0.22%  '
$ws.Range('D38').Value = '''0.0426'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.50%  '
$ws.Range('D39').Value = '''3.14'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.24%  '
$ws.Range('D40').Value = '''2.73'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.87%  '
$ws.Range('D41').Value = '''17.08'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.95%  '
$ws.Range('D42').Value = '''1.88'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.22%  '
$ws.Range('E43').Value = '  -4.95%  '
$ws.Range('D44').Value = '''22.49'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.88%  '
$ws.Range('D45').Value = '''117.68'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.51%  '
$ws.Range('E46').Value = '  -3.81%  '
$ws.Range('D47').Value = '2.056.46'
$ws.Range('E47').Value = '  -3.64%  '
$ws.Range('E48').Value = '  -5.87%  '
$ws.Range('D49').Value = '''3.21'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.18%  '
$ws.Range('D50').Value = '3.209.05'
$ws.Range('E50').Value = '  -0.51%  '
$ws.Range('D51').Value = '''0.238'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.50%  '
